# Auto-generated: applies the Golem_Profits data-refresh diff to all affected sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 174.5
$ws.Range("I2").Value = 174.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 174.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -61.5
$ws.Range("N2").Value = ""

$ws.Range("H4").Value = 350.625
$ws.Range("J4").Value = 211
$ws.Range("L4").Value = 211
$ws.Range("N4").Value = -439

$ws.Range("H40").Value = 1924.138
$ws.Range("I40").Value = 1548.1482
$ws.Range("K40").Value = 1548.1482
$ws.Range("M40").Value = -1373.1482

$ws.Range("H93").Value = 60000
$ws.Range("J93").Value = 60000
$ws.Range("L93").Value = 60000
$ws.Range("N93").Value = -64992

$ws.Range("H113").Value = 2450
$ws.Range("I113").Value = 2450
$ws.Range("K113").Value = 2450
$ws.Range("M113").Value = 804

$ws.Range("H135").Value = 719
$ws.Range("I135").Value = 719
$ws.Range("K135").Value = 6471
$ws.Range("M135").Value = -3936

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 490
$ws.Range("I2").Value = 490
$ws.Range("K2").Value = 490
$ws.Range("M2").Value = -377

$ws.Range("H17").Value = 410.5
$ws.Range("J17").Value = 410.5
$ws.Range("L17").Value = 410.5
$ws.Range("N17").Value = -756.5

$ws.Range("H45").Value = 1326.3334
$ws.Range("I45").Value = 1326.3334
$ws.Range("K45").Value = 1326.3334
$ws.Range("M45").Value = -949.3334

$ws.Range("H61").Value = 2159.5
$ws.Range("I61").Value = 2054.3333
$ws.Range("J61").Value = 2475
$ws.Range("K61").Value = 2054.3333
$ws.Range("L61").Value = 2475
$ws.Range("M61").Value = -1842.3333
$ws.Range("N61").Value = -2899

$ws.Range("H74").Value = 1229.5
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = ""

$ws.Range("H77").Value = 1229.5
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = ""

$ws.Range("H116").Value = 490
$ws.Range("I116").Value = 490
$ws.Range("K116").Value = 490
$ws.Range("M116").Value = 1804

$ws.Range("H132").Value = 3171
$ws.Range("I132").Value = 3999.5
$ws.Range("J132").Value = 1514
$ws.Range("K132").Value = 11998.5
$ws.Range("L132").Value = 4542
$ws.Range("M132").Value = -9468.5
$ws.Range("N132").Value = -9602

$ws.Range("H136").Value = 2159.5
$ws.Range("I136").Value = 2054.3333
$ws.Range("J136").Value = 2475
$ws.Range("K136").Value = 6162.999899999999
$ws.Range("L136").Value = 7425
$ws.Range("M136").Value = -3612.999899999999
$ws.Range("N136").Value = -12525

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 490
$ws.Range("I3").Value = 490
$ws.Range("K3").Value = 490
$ws.Range("M3").Value = -376

$ws.Range("H86").Value = 5273.4546
$ws.Range("I86").Value = 3876.25
$ws.Range("K86").Value = 3876.25
$ws.Range("M86").Value = -2753.25

$ws.Range("H89").Value = 5273.4546
$ws.Range("I89").Value = 3876.25
$ws.Range("K89").Value = 19381.25
$ws.Range("M89").Value = -13765.25

$ws.Range("H107").Value = 1653.125
$ws.Range("I107").Value = 1587.6
$ws.Range("J107").Value = 1980.75
$ws.Range("K107").Value = 1587.6
$ws.Range("L107").Value = 1980.75
$ws.Range("M107").Value = 332.4000000000001
$ws.Range("N107").Value = -5820.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 797.5
$ws.Range("I22").Value = 797.5
$ws.Range("K22").Value = 797.5
$ws.Range("M22").Value = -447.5

$ws.Range("H35").Value = 3887.25
$ws.Range("I35").Value = 3887.25
$ws.Range("K35").Value = 3887.25
$ws.Range("M35").Value = -3593.25

$ws.Range("H59").Value = 65000
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 65000
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 65000
$ws.Range("N59").Value = -67290
$ws.Range("M59").Value = ""

$ws.Range("H62").Value = 2850
$ws.Range("I62").Value = 2850
$ws.Range("K62").Value = 2850
$ws.Range("M62").Value = -2226

$ws.Range("H65").Value = 2850
$ws.Range("I65").Value = 2850
$ws.Range("K65").Value = 14250
$ws.Range("M65").Value = -11130

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 14.454545
$ws.Range("I2").Value = 8.357142
$ws.Range("J2").Value = 25.125
$ws.Range("K2").Value = 50.142852
$ws.Range("L2").Value = 150.75
$ws.Range("M2").Value = 62.857148
$ws.Range("N2").Value = -376.75

$ws.Range("H17").Value = 397.7143
$ws.Range("I17").Value = 217
$ws.Range("J17").Value = 470
$ws.Range("K17").Value = 651
$ws.Range("L17").Value = 1410
$ws.Range("M17").Value = -482
$ws.Range("N17").Value = -1748

$ws.Range("H19").Value = 9000
$ws.Range("J19").Value = 9000
$ws.Range("L19").Value = 27000
$ws.Range("N19").Value = -27348

$ws.Range("H21").Value = 237.5
$ws.Range("J21").Value = 237.5
$ws.Range("L21").Value = 712.5
$ws.Range("N21").Value = -1058.5

$ws.Range("H132").Value = 1766.6666
$ws.Range("J132").Value = 1766.6666
$ws.Range("L132").Value = 15899.9994
$ws.Range("N132").Value = -20959.9994

$ws.Range("H137").Value = 5999.5
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 5999.5
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 17998.5
$ws.Range("N137").Value = -28198.5
$ws.Range("M137").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3312.1667
$ws.Range("I122").Value = 3243.5
$ws.Range("J122").Value = 3449.5
$ws.Range("K122").Value = 9730.5
$ws.Range("L122").Value = 10348.5
$ws.Range("M122").Value = -7280.5
$ws.Range("N122").Value = -15248.5

$ws.Range("H136").Value = 4779.8335
$ws.Range("I136").Value = 4778.8
$ws.Range("K136").Value = 14336.4
$ws.Range("M136").Value = -11786.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 13070
$ws.Range("I51").Value = 13070
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 13070
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -12560
$ws.Range("N51").Value = ""

$ws.Range("H52").Value = 5999.5
$ws.Range("I52").Value = 5999.5
$ws.Range("K52").Value = 5999.5
$ws.Range("M52").Value = -5773.5

$ws.Range("H136").Value = 6110.6665
$ws.Range("I136").Value = 6749.5
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 20248.5
$ws.Range("M136").Value = -17698.5
$ws.Range("N136").Value = -8100
